# Update the credential test data on the active sheet (Sheet3) to the new
# naming scheme, and move the selection from B7 to B9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "dhan0001"
$ws.Range("B2").Value = "dhan#0001"
$ws.Range("A3").Value = "dhan0002"
$ws.Range("B3").Value = "dhan#0002"

$ws.Range("B9").Select()
